{"js": "// Remove the \"Requisitos\" section: the \"Requisitos\" Heading2 paragraph and\n// the following \"LOM3238 -  Projeto Integrado  (Requisito)\" bullet paragraph\n// that lists it as a prerequisite.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text,items/styleBuiltIn\");\nawait context.sync();\n\nconst toDelete = [];\nfor (const p of paragraphs.items) {\n  const text = (p.text || \"\").trim();\n  const isRequisitosHeading =\n    p.styleBuiltIn === Word.BuiltInStyleName.heading2 && text === \"Requisitos\";\n  const isRequisitoBullet = text.indexOf(\"LOM3238\") !== -1;\n\n  if (isRequisitosHeading || isRequisitoBullet) {\n    toDelete.push(p);\n  }\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"Requisitos\" section: the \"Requisitos\" Heading2 paragraph and\n# the following \"LOM3238 -  Projeto Integrado  (Requisito)\" bullet paragraph\n# that lists it as a prerequisite.\n$d = $word.ActiveDocument\n\n# Walk paragraphs back-to-front so deleting one doesn't shift the index of\n# paragraphs we still need to visit.\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $text = $p.Range.Text\n\n    $isRequisitosHeading = ($p.Style.NameLocal -eq \"Heading 2\") -and ($text -match \"^Requisitos\\s*$\")\n    $isRequisitoBullet = $text -match \"LOM3238\"\n\n    if ($isRequisitosHeading -or $isRequisitoBullet) {\n        $p.Range.Delete()\n    }\n}\n"}
